# Weekly update: insert a new price record for "Espinaca" at
# "Feria Lagunitas de Puerto Montt" as the first data row of the block
# (row 58), pushing the existing rows 58-77 down to rows 59-78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 58; this shifts rows 58:77 down to 59:78,
# carrying their styles/formatting along automatically.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new weekly record. Columns that are
# constant across every record in this block are filled in explicitly as
# well so the new row matches the rest of the table.
$ws.Cells.Item(58, 1).Value = 4
$ws.Cells.Item(58, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(58, 3).Value = "Los Lagos"
$ws.Cells.Item(58, 4).Value = 45146
$ws.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(58, 5).Value = 10
$ws.Cells.Item(58, 6).Value = 100112012
$ws.Cells.Item(58, 7).Value = "Espinaca"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 40
$ws.Cells.Item(58, 11).Value = 13500
$ws.Cells.Item(58, 12).Value = 14000
$ws.Cells.Item(58, 13).Value = 13750
$ws.Cells.Item(58, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(58, 15).Value = "Región Metropolitana"
$ws.Cells.Item(58, 16).Value = 1375
$ws.Cells.Item(58, 17).Value = 10
$ws.Cells.Item(58, 18).Value = "Hortaliza"
